$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg).
# This reorders/refreshes the per-row data (a weekly data refresh), while
# columns A,B,C,E,F,G,H,I,J,K,L,Q,R,T remain unchanged.

$rows = @{
    2  = @{ D = 44320; M = 80;  N = 16000; O = 17000; P = 16500; S = 825 }
    3  = @{ D = 44798; M = 80;  N = 21000; O = 22000; P = 21500; S = 1075 }
    4  = @{ D = 44708; M = 80;  N = 20000; O = 21000; P = 20500; S = 1025 }
    5  = @{ D = 44792; M = 100; N = 21000; O = 22000; P = 21500; S = 1075 }
    6  = @{ D = 44893; M = 80;  N = 21000; O = 22000; P = 21625; S = 1081 }
    7  = @{ D = 44761; M = 100; N = 20000; O = 21000; P = 20500; S = 1025 }
    8  = @{ D = 44533; M = 100; N = 16000; O = 17000; P = 16500; S = 825 }
    9  = @{ D = 44890; M = 80;  N = 20000; O = 23000; P = 22250; S = 1112 }
    10 = @{ D = 44357; M = 100; N = 14000; O = 15000; P = 14500; S = 725 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value2  = $vals.D   # column D - Fecha
    $ws.Cells.Item($r, 13).Value2 = $vals.M   # column M - Volumen
    $ws.Cells.Item($r, 14).Value2 = $vals.N   # column N - Precio minimo
    $ws.Cells.Item($r, 15).Value2 = $vals.O   # column O - Precio maximo
    $ws.Cells.Item($r, 16).Value2 = $vals.P   # column P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value2 = $vals.S   # column S - Precio $/Kg
}
